$d = $word.ActiveDocument

# 1. Change "November" to "Desember" in the date line "Bandung, November 2019"
$d.Content.Find.Execute("November", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Desember", 2) | Out-Null

# After the replace, the relevant tail paragraphs (1-indexed via $d.Paragraphs) are:
#   62: " Bandung, Desember 2019"                 -> keep as is
#   63: empty, style BodyTextIndent, right-just.    -> delete entirely
#   64: line of 8 tab runs, right-just., spacing480  -> keep only the first tab run
#   65: empty, spacing480, right-just.               -> delete entirely
#   66: empty, spacing480, right-just.               -> delete entirely
#   67: empty, lang id-ID (document's final para)    -> delete entirely
#
# Paragraph 67 is the very last paragraph in the document body, so a plain
# Range.Delete() on it alone is a no-op (Word always keeps a final paragraph
# mark). Instead: first trim paragraph 64 down to a single tab run, then
# delete the whole range from the end of paragraph 64 through to the end of
# the document in one go - this removes paragraphs 65, 66 and 67 while the
# lone remaining paragraph mark inherits paragraph 64's own formatting.

# Trim paragraph 64's run of tabs down to just the first one.
$p64 = $d.Paragraphs.Item(64)
$tabsStart = $p64.Range.Start
$tabsEnd = $p64.Range.End
$d.Range($tabsStart + 1, $tabsEnd - 1).Delete() | Out-Null

# Remove everything from the (now-trimmed) end of paragraph 64 through the
# end of the document - this deletes paragraphs 65, 66 and 67 in one shot.
$p64 = $d.Paragraphs.Item(64)
$d.Range($p64.Range.End, $d.Content.End).Delete() | Out-Null

# Finally, delete the empty BodyTextIndent paragraph (63) that sits between
# the date line and the tab-run paragraph.
$d.Paragraphs.Item(63).Range.Delete() | Out-Null
